# Append a new data row (85) to the ModCounts sheet:
#   A85 = "2026/02/03"  (text, not an auto-converted date serial)
#   B85 = "逃离鸭科夫"
#   C85 = 1170 (numeric)
# styled the same as the preceding rows (centered alignment, style index 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column A's new cell to stay literal text instead of being
# auto-parsed into a date serial number when it looks like "yyyy/mm/dd".
$ws.Range("A85").NumberFormat = "@"
$ws.Range("A85").Value = "2026/02/03"
$ws.Range("B85").Value = "逃离鸭科夫"
$ws.Range("C85").Value = 1170

# Match the formatting (centered horizontal/vertical alignment) used by
# every other data row by copying row 84's format onto the new row.
$ws.Range("A84:C84").Copy()
$ws.Range("A85:C85").PasteSpecial(-4122)
